# Remove the "Sheet1" worksheet entirely (the formula/demo-sample sheet).
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Sheet1").Delete()

# The remaining data sheet ("Sheet2") becomes the new second sheet; rename it
# to lowercase "sheet1" and trim the extra columns (F1, G2, G4) that are no
# longer part of the table.
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "sheet1"
$ws.Range("F1").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("G4").ClearContents()
